$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MasterTestDataSheet")
$ws2 = $wb.Worksheets.Item("MasterCaseCreation")

# --- Create the new CongaTemplateCreation sheet as a copy of MasterCaseCreation ---
# (keeps the shared header/data styling used on MasterCaseCreation)
$ws2.Copy($null, $ws2) | Out-Null
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "CongaTemplateCreation"

# New App / Name columns on the CongaTemplateCreation sheet
$ws3.Range("D1").Value = "App"
$ws3.Range("D2").Value = "Conga Templates"

# --- Sheet 1: MasterTestDataSheet - new Conga Templates test-case row ---
$ws1.Range("A3").Value = "TC_CongaTemplates_ActualTreatment"
$ws1.Range("B3").Value = "CongaTemplateCreation"
$ws1.Range("C2").Value = "https://test.salesforce.com/"
$ws1.Range("D2").Value = "..\\JunoAutomation\src\resources\\Juno_TestDataSheet.xlsx"

$ws3.Range("E1").Value = "Name"
$ws3.Range("E2").Value = "Auto_TueJan081557572019"

$ws1.Range("C3").Value = "https://test.salesforce.com/"
$ws1.Range("D3").Value = "..\\JunoAutomation\src\resources\\Juno_TestDataSheet.xlsx"
$ws1.Range("E3").Value = "GoogleChrome"

# Restore the quote-prefixed ("as typed text") formatting on row 2 / row 3 cells
# that Value-assignment above would otherwise have cleared.
$ws1.Range("B4").Copy() | Out-Null
$ws1.Range("C2:D2").PasteSpecial(-4122) | Out-Null
$ws1.Range("A3:E3").PasteSpecial(-4122) | Out-Null

# Re-apply the header style to the new D1/E1 cells on CongaTemplateCreation
$ws3.Range("A1").Copy() | Out-Null
$ws3.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$ws3.Range("D1").Value = "App"
$ws3.Range("E1").Value = "Name"

# --- Sheet 2: MasterCaseCreation - blank header-styled D1 cell ---
$ws2.Range("C1").Copy() | Out-Null
$ws2.Range("D1").PasteSpecial(-4122) | Out-Null
$ws2.Range("D1").Value = ""

# --- Selections / active sheet ---
$ws2.Activate() | Out-Null
$ws2.Range("C6").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("D3").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("F7").Select() | Out-Null
